# Swap the "group-name" and "group-code" columns (C and D) on the active sheet.
# Before: A=code, B=status, C=codeforiati:group-name, D=codeforiati:group-code
# After:  A=code, B=status, C=codeforiati:group-code, D=codeforiati:group-name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
